$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Inscricoes")

# Row 27
$ws.Range("E27").Value = 17

# Row 40
$ws.Range("E40").Value = 27
$ws.Range("G40").Value = 2
$ws.Range("H40").Value = 21

# Row 42
$ws.Range("E42").Value = 41
$ws.Range("F42").Value = 25
$ws.Range("H42").Value = 34

# Row 43
$ws.Range("E43").Value = 32
$ws.Range("F43").Value = 17
$ws.Range("H43").Value = 20

# Row 44
$ws.Range("E44").Value = 31

# Row 50
$ws.Range("E50").Value = 32
$ws.Range("G50").Value = 9
$ws.Range("H50").Value = 21
